$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.178.75"
$ws.Range("E2").Value = "  -4.00%  "
$ws.Range("D3").Value = "'1.657.31"
$ws.Range("E3").Value = "  -2.65%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "'217.77"
$ws.Range("E5").Value = "  -2.50%  "
$ws.Range("D6").Value = "'0.5143"
$ws.Range("E6").Value = "  -3.05%  "
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("D8").Value = "'0.2580"
$ws.Range("E8").Value = "  -2.91%  "
$ws.Range("D9").Value = "'0.06433"
$ws.Range("E9").Value = "  -2.17%  "
$ws.Range("D10").Value = "'19.94"
$ws.Range("D11").Value = "'0.07806"
$ws.Range("E11").Value = "  +2.45%  "
$ws.Range("D12").Value = "'1.659.63"
$ws.Range("E12").Value = "  -2.82%  "
$ws.Range("D13").Value = "'4.295"
$ws.Range("D14").Value = "'1.885.76"
$ws.Range("E14").Value = "  -2.64%  "
$ws.Range("D15").Value = "'0.5545"
$ws.Range("E15").Value = "  -3.93%  "
$ws.Range("D16").Value = "'0.0₅8057"
$ws.Range("E16").Value = "  -1.05%  "
$ws.Range("D17").Value = "'64.20"
$ws.Range("D18").Value = "'26.200.78"
$ws.Range("E18").Value = "  -3.92%  "
$ws.Range("D19").Value = "'1.005"
$ws.Range("E19").Value = "  +0.21%  "
$ws.Range("D20").Value = "'211.21"
$ws.Range("E20").Value = "  -1.83%  "
$ws.Range("D21").Value = "'4.423"
$ws.Range("E21").Value = "  -3.95%  "
$ws.Range("E22").Value = "  -2.89%  "
$ws.Range("D23").Value = "'5.973"
$ws.Range("E23").Value = "  +0.25%  "
$ws.Range("D24").Value = "'1.005"
$ws.Range("E24").Value = "  +0.26%  "
$ws.Range("D25").Value = "'143.47"
$ws.Range("D26").Value = "'1.755"
$ws.Range("E26").Value = "  +3.23%  "
$ws.Range("E27").Value = "  -2.73%  "
$ws.Range("D28").Value = "'6.969"
$ws.Range("E28").Value = "  -3.14%  "
$ws.Range("D29").Value = "'15.77"
$ws.Range("E29").Value = "  -2.15%  "
$ws.Range("D30").Value = "'0.05212"
$ws.Range("E30").Value = "  -2.67%  "
$ws.Range("D31").Value = "'1.253"
$ws.Range("E31").Value = "  -2.44%  "
$ws.Range("D32").Value = "'3.365"
$ws.Range("E32").Value = "  -2.64%  "
$ws.Range("D33").Value = "'3.216"
$ws.Range("E33").Value = "  -5.25%  "
$ws.Range("D34").Value = "'1.567"
$ws.Range("E34").Value = "  -4.35%  "
$ws.Range("D35").Value = "'2.760"
$ws.Range("E35").Value = "  -3.66%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "'2.372"
$ws.Range("E36").Value = "  -1.62%  "
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").Value = "'0.9301"
$ws.Range("E37").Value = "  -1.54%  "
$ws.Range("D38").Value = "'1.172.27"
$ws.Range("E38").Value = "  +12.81%  "
$ws.Range("D39").Value = "'0.5693"
$ws.Range("E39").Value = "  -1.79%  "
$ws.Range("E40").Value = "  -1.91%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "'0.8449"
$ws.Range("E41").Value = "  +0.74%  "
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").Value = "'1.005"
$ws.Range("E42").Value = "  +0.26%  "
$ws.Range("D43").Value = "'5.668"
$ws.Range("E43").Value = "  -1.87%  "
$ws.Range("D44").Value = "'100.53"
$ws.Range("E44").Value = "  -0.29%  "
$ws.Range("E45").Value = "  -2.68%  "
$ws.Range("E46").Value = "  +1.01%  "
$ws.Range("E47").Value = "  +0.41%  "
$ws.Range("D48").Value = "'55.89"
$ws.Range("E48").Value = "  -3.13%  "
$ws.Range("E49").Value = "  -0.13%  "
$ws.Range("D50").Value = "'7.900"
$ws.Range("E50").Value = "  -1.59%  "
$ws.Range("D51").Value = "'0.05055"
$ws.Range("E51").Value = "  -3.28%  "
